$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("H1").Value = "Roth conv"
    $ws.Range("C1").Value = "taxable ctrb"
    $ws.Range("D1").Value = "401k ctrb"
    $ws.Range("E1").Value = "Roth 401k ctrb"
    $ws.Range("F1").Value = "IRA ctrb"
    $ws.Range("G1").Value = "Roth IRA ctrb"

    $ws.Range("A1:I1").Font.Bold = $true
}

$null = $ws2.Activate()
$null = $ws2.Range("A1:XFD1").Select()

$null = $ws1.Activate()
$null = $ws1.Range("A1:XFD1").Select()
